$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, copying the format of the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value2 = "Save"

# Fill H2:H61 - "Save" indicator: 1 if the row's sum (column G) exceeds 9, else 0
for ($r = 2; $r -le 61; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -gt 9) {
        $ws.Cells.Item($r, 8).Value2 = 1
    } else {
        $ws.Cells.Item($r, 8).Value2 = 0
    }
}
